$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.49"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'20.95"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").Value = "'0.06197"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").Value = "'3.581"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'6.548"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'1.482"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.8227"
$ws.Range("D9").ClearFormats()
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1642"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08224"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03510"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03117"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09126"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.769"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001615"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04695"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006416"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "17TigerCashTCHBestin24h"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.006144"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001068"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001502"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.822"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.320"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01386"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "23OneONE"
$ws.Range("D25").Value = "'0.3302"
$ws.Range("D25").ClearFormats()
$ws.Range("D40").Value = "'0.04675"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "'0.007018"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").Value = "'0.004605"
$ws.Range("D42").ClearFormats()
$ws.Range("D45").Value = "'0.00006301"
$ws.Range("D45").ClearFormats()
$ws.Range("D47").Value = "'0.8461"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = "'0.001970"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.00001902"
$ws.Range("D49").ClearFormats()
